$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Add the required "Experimental" element (row 7: Experimental / true).
# A plain `.Value = "true"` gets auto-coerced to a native Excel boolean
# (stored as t="b"), but the target workbook stores the literal text "true"
# as a shared string (t="s"), matching its existing column-B cells. Entering
# it as a quoted-text formula and then doing the classic "convert formula to
# value" (Copy + Paste Special -> Values) keeps it as literal text without
# touching the cell's style.
$cell = $ws.Cells.Item(7, 2)
$cell.Formula = "=""true"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Update the Date value (row 8: Date / <timestamp>)
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
